$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "UVA Digital"
$ws.Range("B2").Value = " Skill set: Well versed & trained with latest versions of design software s like: 3D max..."
$ws.Range("C2").Value = "1-2 Yrs"
$ws.Range("D2").Value = "Bangalore/Bengaluru"
$ws.Range("E2").Value = "Not disclosed"
$ws.Range("F2").Value = "['3D max', '3D Design', 'Maya', 'Photoshop', 'CorelDraw', 'Visual Designer', '2D Design']"

$ws.Range("A3").Value = "Acme Designers"
$ws.Range("B3").Value = " You must have experience in AutoCAD, 3D Studio Max or Cinema 4D"
$ws.Range("C3").Value = "2-7 Yrs"
$ws.Range("D3").Value = "Bangalore/Bengaluru(HSR Layout)"
$ws.Range("E3").Value = "Not disclosed"
$ws.Range("F3").Value = "['Visualiser', 'Lumion', 'Sketchup', '3D Design', 'Illustrator', 'AutoCAD', 'Google Sketchup', 'REVIT']"

$ws.Range("A4").Value = "NEWAGE PRODUCT DESIGNS PRIVATE LIMITED"
$ws.Range("B4").Value = " Selected candidates will be required to complete a Design Test Package.Bachelors Degree..."
$ws.Range("C4").Value = "4-7 Yrs"
$ws.Range("D4").Value = "Mumbai"
$ws.Range("E4").Value = "Not disclosed"
$ws.Range("F4").Value = "['Design', 'Execution', 'Project Management', 'Rendering', 'Lighting', 'Texturing', 'Interiors', 'New Product']"

$ws.Range("A5").Value = "Lightcraft"
$ws.Range("B5").Value = " Candidates with their own vehicle preferred2+ Years Experience"
$ws.Range("C5").Value = "2-7 Yrs"
$ws.Range("D5").Value = "New Delhi(Okhla)"
$ws.Range("E5").Value = "Not disclosed"
$ws.Range("F5").Value = "['Design', '3D Modeling', 'Project Coordination', 'Lead Generation', 'Project Sales', 'Bdm', 'Business Development Management']"

$ws.Range("A6").Value = "Yazaki India Private Limited"
$ws.Range("B6").Value = " 3D master and 2D master release in teamcenter for each harness family CCDImplement CCD ..."
$ws.Range("C6").Value = "4-8 Yrs"
$ws.Range("D6").Value = "Chennai(Perungudi)"
$ws.Range("E6").Value = "Not disclosed"
$ws.Range("F6").Value = "['Design', 'VAVE', '3D Cad', 'UG NX', 'Teamcenter']"

$ws.Range("A7").Value = "Gokaldas Exports Ltd"
$ws.Range("B7").Value = " To have 3 to 5 years of experience in handling 3D soft wear preferably BROWZWEAR / CLO"
$ws.Range("C7").Value = "3-5 Yrs"
$ws.Range("D7").Value = "Bangalore/Bengaluru"
$ws.Range("E7").Value = "Not disclosed"
$ws.Range("F7").Value = "['Sewing', 'CAD', 'Photoshop']"

$ws.Range("A8").Value = "Toppr"
$ws.Range("B8").Value = " Certificate, associates degree, or bachelors degree in graphic design or a related fiel..."
$ws.Range("C8").Value = "3-7 Yrs"
$ws.Range("D8").Value = "Hyderabad/Secunderabad"
$ws.Range("E8").Value = "Not disclosed"
$ws.Range("F8").Value = "['Graphics', 'Visual Effects', '3D', 'Kaizen', 'Time management', 'Texturing', 'Venture capital', 'Maya']"

$ws.Range("A9").Value = "Gokaldas Exports Ltd"
$ws.Range("B9").Value = " To have 3 to 5 years of experience in handling 3D soft wear preferably BROWZWEAR / CLO"
$ws.Range("C9").Value = "3-5 Yrs"
$ws.Range("D9").Value = "Bangalore/Bengaluru"
$ws.Range("E9").Value = "Not disclosed"
$ws.Range("F9").Value = "['Sewing', 'CAD', 'Photoshop']"

$ws.Range("A10").Value = "The Patina Studio"
$ws.Range("B10").Value = " Must be able to read architectural drawings from AutoCAD, with knowledge of basic inter..."
$ws.Range("C10").Value = "2-6 Yrs"
$ws.Range("D10").Value = "Delhi / NCR"
$ws.Range("E10").Value = "3,50,000 - 6,00,000 PA."
$ws.Range("F10").Value = "['Sketchup', 'Visualiser', 'Architecture', 'AutoCAD', 'Interiors', 'Photoshop', 'VRAY', '3Ds Max']"

$ws.Range("A11").Value = "Vedarth Animation Studio Pvt. Ltd."
$ws.Range("B11").Value = " Minimum 1+ years of working knowledge of and production experience with After EffectsNO..."
$ws.Range("C11").Value = "1-3 Yrs"
$ws.Range("D11").Value = "Pune"
$ws.Range("E11").Value = "Not disclosed"
$ws.Range("F11").Value = "['Visual Effects', 'Adobe Premiere Pro', 'VFX', 'Adobe After Effects', 'Compositing', 'VRAY', '3Ds Max', '3D Compositing']"

$ws.Range("A12").Value = "Quest Global"
$ws.Range("B12").Value = " We are looking for a 2D/3D Graphic Artist with minimum experience of 4-6 years to work ..."
$ws.Range("C12").Value = "3-6 Yrs"
$ws.Range("D12").Value = "Pune"
$ws.Range("E12").Value = "Not disclosed"
$ws.Range("F12").Value = "['3D Graphics', 'Design', 'UX', 'Illustrator', 'Maya', '3D Animation', 'Photoshop', 'Heavy Engineering']"

$ws.Range("A13").Value = "Advids.co"
$ws.Range("B13").Value = " Why you should DEFINITELY apply & join us Excellent communication skills for 3D Particl..."
$ws.Range("C13").Value = "1-5 Yrs"
$ws.Range("D13").Value = "Ahmedabad, Jaipur, Surat"
$ws.Range("E13").Value = "2,00,000 - 6,00,000 PA."
$ws.Range("F13").Value = "['Simulation Artist', 'Houdini', 'FX Artist', 'Fume FX', 'Time Management', 'Technical Skills', '3Ds Max', 'Particle Simulation']"

$ws.Range("A14").Value = "Advids.co"
$ws.Range("B14").Value = " Why you should DEFINITELY apply & join us Excellent communication skills for 3D Particl..."
$ws.Range("C14").Value = "1-5 Yrs"
$ws.Range("D14").Value = "Chandigarh, Lucknow, Delhi / NCR"
$ws.Range("E14").Value = "2,00,000 - 6,00,000 PA."
$ws.Range("F14").Value = "['Simulation Artist', 'Houdini', 'FX Artist', 'Fume FX', 'Time Management', 'Technical Skills', '3Ds Max', 'Particle Simulation']"

$ws.Range("A15").Value = "Advids.co"
$ws.Range("B15").Value = " Why you should DEFINITELY apply & join us Excellent communication skills for 3D Particl..."
$ws.Range("C15").Value = "1-5 Yrs"
$ws.Range("D15").Value = "Guwahati, Bhubaneswar, Kolkata"
$ws.Range("E15").Value = "2,00,000 - 6,00,000 PA."
$ws.Range("F15").Value = "['Simulation Artist', 'Houdini', 'FX Artist', 'Fume FX', 'Time Management', 'Technical Skills', '3Ds Max', 'Particle Simulation']"

$ws.Range("A16").Value = "Advids.co"
$ws.Range("B16").Value = " He / she should have good handle on various style of art for 3D Motion graphics video c..."
$ws.Range("C16").Value = "1-5 Yrs"
$ws.Range("D16").Value = "Madurai, Chennai, Coimbatore"
$ws.Range("E16").Value = "2,00,000 - 6,00,000 PA."
$ws.Range("F16").Value = "['Cinema 4D', '3D Character Animation', 'Maya', 'Art', 'Unity3D', '3D Maya', '3Ds Max', '3D Compositing']"

$ws.Range("A17").Value = "Advids.co"
$ws.Range("B17").Value = " He / she should have good handle on various style of art for 3D Motion graphics video c..."
$ws.Range("C17").Value = "1-5 Yrs"
$ws.Range("D17").Value = "Bhopal, Gwalior, Indore"
$ws.Range("E17").Value = "2,00,000 - 6,00,000 PA."
$ws.Range("F17").Value = "['Cinema 4D', '3D Character Animation', 'Maya', 'Art', 'Unity3D', '3D Maya', '3Ds Max', '3D Compositing']"

$ws.Range("A18").Value = "Advids.co"
$ws.Range("B18").Value = " He / she should have a good handle on various style of art for 3D Motion graphics video..."
$ws.Range("C18").Value = "1-5 Yrs"
$ws.Range("D18").Value = "Vijayawada, Visakhapatnam, Hyderabad/Secunderabad"
$ws.Range("E18").Value = "2,00,000 - 6,00,000 PA."
$ws.Range("F18").Value = "['Cinema 4D', '3D Generalist', 'Maya', 'texturing', '3D Modeler', 'Mudbox', 'Rigging', '3Ds Max']"

$ws.Range("A19").Value = "Advids.co"
$ws.Range("B19").Value = " He / she should have a good handle on various style of art for 3D Motion graphics video..."
$ws.Range("C19").Value = "1-5 Yrs"
$ws.Range("D19").Value = "Mysore/Mysuru, Bangalore/Bengaluru, Belagavi/Belgaum"
$ws.Range("E19").Value = "2,00,000 - 6,00,000 PA."
$ws.Range("F19").Value = "['Cinema 4D', '3D Generalist', 'Maya', 'texturing', '3D Modeler', 'Mudbox', 'Rigging', '3Ds Max']"

$ws.Range("A20").Value = "Advids.co"
$ws.Range("B20").Value = " He / she should have a good handle on various style of art for 3D Motion graphics video..."
$ws.Range("C20").Value = "1-5 Yrs"
$ws.Range("D20").Value = "Kochi/Cochin, Vellore, Trivandrum/Thiruvananthapuram"
$ws.Range("E20").Value = "2,00,000 - 6,00,000 PA."
$ws.Range("F20").Value = "['Cinema 4D', '3D Generalist', 'Maya', 'texturing', '3D Modeler', 'Mudbox', 'Rigging', '3Ds Max']"

$ws.Range("A21").Value = "UVA Digital"
$ws.Range("B21").Value = " Skill set: Well versed & trained with latest versions of design software s like: 3D max..."
$ws.Range("C21").Value = "1-2 Yrs"
$ws.Range("D21").Value = "Bangalore/Bengaluru"
$ws.Range("E21").Value = "Not disclosed"
$ws.Range("F21").Value = "['3D max', '3D Design', 'Maya', 'Photoshop', 'CorelDraw', 'Visual Designer', '2D Design']"
